$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 46: 四方坪站 (station index 4 -> shared string "四方坪站")
$ws.Cells.Item(46, 1).Value = 45923
$ws.Cells.Item(46, 2).Value = "四方坪站"
$ws.Cells.Item(46, 3).Value = 8507.1
$ws.Cells.Item(46, 4).Value = 6928.54
$ws.Cells.Item(46, 5).Value = 2978.9
$ws.Cells.Item(46, 6).Value = 378

# New row 47: 高岭站 (station index 5 -> shared string "高岭站")
$ws.Cells.Item(47, 1).Value = 45923
$ws.Cells.Item(47, 2).Value = "高岭站"
$ws.Cells.Item(47, 3).Value = 5170.99
$ws.Cells.Item(47, 4).Value = 4034.1
$ws.Cells.Item(47, 5).Value = 1387.38
$ws.Cells.Item(47, 6).Value = 189

# Copy styles from the row above (row 45) to keep number formats consistent
$ws.Range("A45:F45").Copy() | Out-Null
$ws.Range("A46:F46").PasteSpecial(-4122) | Out-Null
$ws.Range("A47:F47").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply values in case paste-special of formats affected anything
$ws.Cells.Item(46, 1).Value = 45923
$ws.Cells.Item(46, 2).Value = "四方坪站"
$ws.Cells.Item(46, 3).Value = 8507.1
$ws.Cells.Item(46, 4).Value = 6928.54
$ws.Cells.Item(46, 5).Value = 2978.9
$ws.Cells.Item(46, 6).Value = 378

$ws.Cells.Item(47, 1).Value = 45923
$ws.Cells.Item(47, 2).Value = "高岭站"
$ws.Cells.Item(47, 3).Value = 5170.99
$ws.Cells.Item(47, 4).Value = 4034.1
$ws.Cells.Item(47, 5).Value = 1387.38
$ws.Cells.Item(47, 6).Value = 189
